# Apply updated cryptocurrency market data (prices and 1h volume change)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.000", "329.70").
# Force it to Text format first so Excel does not coerce it into a Double
# and strip significant trailing zeros / thousand-dot grouping.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.211.91"
$ws.Range("E2").Value = "  +2.85%  "

$ws.Range("D3").Value = "1.815.72"
$ws.Range("E3").Value = "  +4.41%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "329.70"
$ws.Range("E5").Value = "  +1.67%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").Value = "0.4425"
$ws.Range("E7").Value = "  +4.35%  "

$ws.Range("D8").Value = "0.3696"

$ws.Range("D9").Value = "44.61"
$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("D10").Value = "0.07695"
$ws.Range("E10").Value = "  +3.08%  "

$ws.Range("D11").Value = "1.125"
$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").Value = "22.03"
$ws.Range("E13").Value = "  +2.05%  "

$ws.Range("D14").Value = "6.255"
$ws.Range("E14").Value = "  +2.90%  "

$ws.Range("D15").Value = "7.542"
$ws.Range("E15").Value = "  +5.21%  "

$ws.Range("D16").Value = "1.822.06"
$ws.Range("E16").Value = "  +4.95%  "

$ws.Range("D17").Value = "92.77"
$ws.Range("E17").Value = "  +6.61%  "

$ws.Range("E18").Value = "  +1.54%  "

$ws.Range("D19").Value = "0.06601"
$ws.Range("E19").Value = "  +9.10%  "

$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("E21").Value = "  +4.04%  "

$ws.Range("D22").Value = "6.192"
$ws.Range("E22").Value = "  +2.07%  "

$ws.Range("D23").Value = "28.271.34"
$ws.Range("E23").Value = "  +3.01%  "

$ws.Range("D24").Value = "11.67"
$ws.Range("E24").Value = "  +3.10%  "

$ws.Range("D25").Value = "2.024"
$ws.Range("E25").Value = "  -15.97%  "

$ws.Range("D26").Value = "20.75"
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("D27").Value = "155.58"
$ws.Range("E27").Value = "  +4.08%  "

$ws.Range("D28").Value = "2.023.41"
$ws.Range("E28").Value = "  +4.55%  "

$ws.Range("D29").Value = "2.316"
$ws.Range("E29").Value = "  -3.45%  "

$ws.Range("D30").Value = "128.12"
$ws.Range("E30").Value = "  +1.20%  "

$ws.Range("D31").Value = "1.202"
$ws.Range("E31").Value = "  -5.64%  "

$ws.Range("D32").Value = "5.866"
$ws.Range("E32").Value = "  +4.56%  "

$ws.Range("D33").Value = "0.09205"
$ws.Range("E33").Value = "  +1.72%  "

$ws.Range("D34").Value = "3.662"
$ws.Range("E34").Value = "  -1.91%  "

$ws.Range("D35").Value = "13.05"
$ws.Range("E35").Value = "  +4.31%  "

$ws.Range("D36").Value = "0.02351"
$ws.Range("E36").Value = "  +3.01%  "

$ws.Range("D37").Value = "0.2170"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "0.6569"
$ws.Range("E38").Value = "  +2.34%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.06206"
$ws.Range("E39").Value = "  +0.69%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "5.147"
$ws.Range("E40").Value = "  +2.24%  "

$ws.Range("D41").Value = "1.195"
$ws.Range("E41").Value = "  +0.68%  "

$ws.Range("D42").Value = "8.141"
$ws.Range("E42").Value = "  +3.56%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.386"
$ws.Range("E44").Value = "  -2.01%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.84"
$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("D46").Value = "0.6075"
$ws.Range("E46").Value = "  +3.76%  "

$ws.Range("D47").Value = "3.764"
$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("D48").Value = "127.25"
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("D49").Value = "2.034"
$ws.Range("E49").Value = "  +4.90%  "

$ws.Range("E50").Value = "  +5.13%  "

$ws.Range("D51").Value = "0.06977"
$ws.Range("E51").Value = "  +2.22%  "

# Restore default (unstyled) cell formatting for the price column so the
# workbook style table matches its original shape.
$priceRange.Style = "Normal"